$d = $word.ActiveDocument

$replacements = @(
    @{old="91×35=3185"; new="46×23=1058"},
    @{old="62×74=4588"; new="93×18=1674"},
    @{old="68×42=2856"; new="81×39=3159"},
    @{old="67×47=3149"; new="91×54=4914"},
    @{old="26×18=468";  new="21×68=1428"},
    @{old="24×48=1152"; new="45×23=1035"},
    @{old="74×63=4662"; new="38×77=2926"},
    @{old="55×59=3245"; new="14×55=770"},
    @{old="97×90=8730"; new="39×65=2535"},
    @{old="79×80=6320"; new="33×48=1584"},
    @{old="81×67=5427"; new="75×94=7050"},
    @{old="74×43=3182"; new="22×46=1012"},
    @{old="75×57=4275"; new="97×84=8148"},
    @{old="82×30=2460"; new="16×64=1024"},
    @{old="27×55=1485"; new="51×44=2244"},
    @{old="99×54=5346"; new="34×86=2924"},
    @{old="98×35=3430"; new="94×52=4888"},
    @{old="90×79=7110"; new="13×58=754"},
    @{old="30×27=810";  new="72×30=2160"},
    @{old="69×38=2622"; new="81×89=7209"},
    @{old="15×39=585";  new="75×21=1575"},
    @{old="92×92=8464"; new="78×96=7488"},
    @{old="40×13=520";  new="27×48=1296"},
    @{old="33×32=1056"; new="14×21=294"},
    @{old="29×51=1479"; new="42×18=756"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
